# Apply the authored edit to Workbook.xlsx
#
# Summary of changes:
#  - "prioriy List" sheet: fill in "Assigned to" (D) / "Status" (E) columns for
#    the first four user stories (rows 2-5), give row 6 a custom height, and
#    move the sheet's remembered selection to E5.
#  - "Queries" sheet: move that sheet's remembered selection to B18 (without
#    changing which sheet tab is active - "prioriy List" stays active).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("prioriy List")
$ws3 = $wb.Worksheets.Item("Queries")

# --- 1. Cell values -------------------------------------------------------
# Set these in this exact order so the shared-string table grows with the
# same new entries, in the same order, as the authored edit:
#   "Researching", "Carl", "Dimitar & Ciaran"
$ws1.Range("E2").Value = "Researching"
$ws1.Range("D2").Value = "Carl"
$ws1.Range("E3").Value = "Researching"
$ws1.Range("D3").Value = "Carl"
$ws1.Range("E4").Value = "Researching"
$ws1.Range("D4").Value = "Dimitar & Ciaran"
$ws1.Range("E5").Value = "Researching"
$ws1.Range("D5").Value = "Sahar"

# --- 2. Cell formatting (indent levels) -----------------------------------
# Applied after all values so the new cell styles are created in the same
# order as the authored edit (no-indent, then indent=1, then indent=7).
$ws1.Range("E2").IndentLevel = 0
$ws1.Range("D4").IndentLevel = 1
$ws1.Range("D2").IndentLevel = 7

$ws1.Range("E3").IndentLevel = 0
$ws1.Range("E4").IndentLevel = 0
$ws1.Range("E5").IndentLevel = 0
$ws1.Range("D3").IndentLevel = 7
$ws1.Range("D5").IndentLevel = 7

# --- 3. Row height ----------------------------------------------------------
$ws1.Rows.Item(6).RowHeight = 29

# --- 4. Selections ----------------------------------------------------------
# Update the "Queries" sheet's remembered selection first, then return to
# "prioriy List" so it remains the active/visible tab, matching the target.
$ws3.Range("B18").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("E5").Select() | Out-Null
